$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column holds text-formatted numbers (e.g. "30.738.15", "0.9991", "19.70")
# that must stay literal text, not be coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Bulk update Price (D) and Volume(1h) (E) columns for rows with both changes
$ws.Range("D2").Value = "30.738.15"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "1.882.72"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "238.16"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.4751"
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("D8").Value = "0.2822"
$ws.Range("E8").Value = "  +3.26%  "
$ws.Range("D9").Value = "0.06517"
$ws.Range("E9").Value = "  +3.65%  "
$ws.Range("D10").Value = "18.61"
$ws.Range("E10").Value = "  +14.35%  "
$ws.Range("D11").Value = "1.882.71"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("D12").Value = "0.07558"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").Value = "95.12"
$ws.Range("E13").Value = "  +13.37%  "
$ws.Range("D14").Value = "5.074"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").Value = "0.6488"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").Value = "301.92"
$ws.Range("E16").Value = "  +32.01%  "
$ws.Range("D17").Value = "30.723.33"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D20").Value = "0.000007542"
$ws.Range("E20").Value = "  +3.37%  "
$ws.Range("D21").Value = "2.130.20"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").Value = "0.9993"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "5.132"
$ws.Range("E23").Value = "  +4.25%  "
$ws.Range("D24").Value = "6.146"
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("D25").Value = "169.18"
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("D26").Value = "9.231"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "19.70"
$ws.Range("E27").Value = "  +10.57%  "
$ws.Range("D28").Value = "1.941"
$ws.Range("E28").Value = "  +3.52%  "
$ws.Range("D29").Value = "0.1063"
$ws.Range("E29").Value = "  +3.62%  "
$ws.Range("D30").Value = "1.350"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D32").Value = "3.940"
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("D33").Value = "0.05059"
$ws.Range("E33").Value = "  +4.28%  "
$ws.Range("D35").Value = "0.7203"
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("D36").Value = "2.718"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").Value = "0.01916"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D39").Value = "2.051"
$ws.Range("E39").Value = "  +6.53%  "
$ws.Range("D40").Value = "0.8975"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "107.34"
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("D42").Value = "0.9997"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").Value = "0.4183"
$ws.Range("E43").Value = "  +3.73%  "
$ws.Range("D44").Value = "5.596"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("D45").Value = "7.301"
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("D46").Value = "64.80"
$ws.Range("E46").Value = "  +4.44%  "
$ws.Range("D47").Value = "0.1218"
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("D48").Value = "8.955"
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("D49").Value = "34.56"
$ws.Range("E49").Value = "  +4.05%  "
$ws.Range("D50").Value = "0.05592"
$ws.Range("E50").Value = "  +1.48%  "

# Rows where only Volume(1h) changes
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E38").Value = "  +1.82%  "

# Row 18/19 content swap: Dai and Avalanche swap positions with updated values
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "13.06"
$ws.Range("E18").Value = "  +5.51%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  +0.05%  "

# Row 51: Decentraland replaced with NEARProtocol
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.374"
$ws.Range("E51").Value = "  +1.45%  "
